$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 2 - "H" row)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 267
$wsOff.Range("C2").Value = 194
$wsOff.Range("D2").Value = 66
$wsOff.Range("E2").Value = 34

# Update DEF sheet (row 2 - "H" row)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 250
$wsDef.Range("C2").Value = 175
$wsDef.Range("D2").Value = 61
$wsDef.Range("E2").Value = 27
$wsDef.Range("F2").Value = 5
$wsDef.Range("G2").Value = 5

$wb.Save()
